$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: F3 4 -> 5, H3 4 -> 5 (E3 unchanged)
$ws.Range("F3").Value = 5
$ws.Range("H3").Value = 5

# Row 7: E7 14 -> 15
$ws.Range("E7").Value = 15

# Row 8: E8 24 -> 25
$ws.Range("E8").Value = 25

# Row 12: E12 14 -> 15, F12 3 -> 4, H12 3 -> 4
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 4
$ws.Range("H12").Value = 4

# Row 16: E16 208 -> 210, F16 49 -> 51, H16 49 -> 51
$ws.Range("E16").Value = 210
$ws.Range("F16").Value = 51
$ws.Range("H16").Value = 51

# Row 18: E18 52 -> 53
$ws.Range("E18").Value = 53

$wb.Save()
